$d = $word.ActiveDocument

# 1. Timestamp in the first (summary) table - text is unique in the document.
$d.Content.Find.Execute(
    "2023-02-25 02:23:55.360784", $true, $false, $false, $false, $false,
    $true, 1, $false, "2023-02-25 09:45:54.287682", 2)

# The QUESTIONS / ANSWER table is the second table in the document.
$t = $d.Tables(2)

# 2. "Have you been capable of finding humor..." answer (row 2, answer column).
$t.Cell(2, 2).Range.Text = "(2) Definitely less than I used to"

# 3. "Have you anticipated things with pleasure and excitement?" answer (row 3).
$t.Cell(3, 2).Range.Text = "(3) Yes, most of the time"

# 4. "Have you been so unhappy that you have experienced trouble sleeping?" (row 8).
#    NOTE: the text "(2) Yes, quite often" also occurs, unchanged, in row 9, so
#    this must be targeted at the specific cell rather than a document-wide find.
$t.Cell(8, 2).Range.Text = "(1) Not very often"

# 5. "Have you been so unhappy that you have shed tears?" answer (row 10).
$t.Cell(10, 2).Range.Text = "(2) Yes, quite often"

# 6. EPDS Score value (row 11) - preserve the trailing space as in the original.
$t.Cell(11, 2).Range.Text = "(19) "

Write-Host "done"
